{"js": "// The \"MBS Datasets\" Heading-3 paragraph loses its paragraph style / spacing\n// and instead carries the equivalent look (font, color, size) as direct\n// formatting on the paragraph mark, while the trailing blank paragraphs and\n// the paragraph that used to hold the page-break run are collapsed into it\n// (the page-break run itself is kept, now ending the \"MBS Datasets\"\n// paragraph).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the first \"MBS Datasets\" Heading 3 paragraph that is immediately\n// followed by two empty paragraphs and a page-break paragraph (the second\n// \"MBS Datasets\" heading later in the document must stay untouched).\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length - 3; i++) {\n  if (\n    paragraphs.items[i].text === \"MBS Datasets\" &&\n    paragraphs.items[i + 1].text === \"\" &&\n    paragraphs.items[i + 2].text === \"\" &&\n    paragraphs.items[i + 3].text === \"\\f\"\n  ) {\n    targetIndex = i;\n    break;\n  }\n}\n\nif (targetIndex === -1) {\n  throw new Error('Could not locate the target \"MBS Datasets\" heading paragraph.');\n}\n\nconst headingParagraph = paragraphs.items[targetIndex];\nconst blankParagraph1 = paragraphs.items[targetIndex + 1];\nconst blankParagraph2 = paragraphs.items[targetIndex + 2];\nconst pageBreakParagraph = paragraphs.items[targetIndex + 3];\n\n// Replace the heading paragraph's OOXML: drop pStyle/spacing, bake the\n// Heading-3 look into the paragraph mark's rPr, keep the two existing runs\n// untouched, and append the page-break run that used to live in its own\n// paragraph.\nconst newParagraphOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  '<w:p w14:paraId=\"708C0CEA\" w14:textId=\"38B3AB9C\" w:rsidR=\"007325B7\" w:rsidRPr=\"000E1730\" w:rsidRDefault=\"007325B7\" w:rsidP=\"007325B7\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n  \"<w:pPr>\" +\n  \"<w:rPr>\" +\n  '<w:rFonts w:eastAsiaTheme=\"majorEastAsia\" w:cstheme=\"majorBidi\"/>' +\n  '<w:color w:val=\"0F4761\" w:themeColor=\"accent1\" w:themeShade=\"BF\"/>' +\n  '<w:sz w:val=\"28\"/>' +\n  '<w:szCs w:val=\"28\"/>' +\n  '<w:u w:val=\"single\"/>' +\n  \"</w:rPr>\" +\n  \"</w:pPr>\" +\n  \"<w:r>\" +\n  \"<w:rPr>\" +\n  '<w:u w:val=\"single\"/>' +\n  \"</w:rPr>\" +\n  \"<w:t>MBS</w:t>\" +\n  \"</w:r>\" +\n  '<w:r w:rsidRPr=\"000E1730\">' +\n  \"<w:rPr>\" +\n  '<w:u w:val=\"single\"/>' +\n  \"</w:rPr>\" +\n  '<w:t xml:space=\"preserve\"> Datasets</w:t>' +\n  \"</w:r>\" +\n  \"<w:r>\" +\n  \"<w:rPr>\" +\n  \"<w:b/>\" +\n  \"<w:bCs/>\" +\n  \"</w:rPr>\" +\n  '<w:br w:type=\"page\"/>' +\n  \"</w:r>\" +\n  \"</w:p>\" +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\nheadingParagraph.getRange().insertOoxml(newParagraphOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// Remove the two blank paragraphs and the now-redundant paragraph that used\n// to hold only the page-break run (delete from the bottom up so the other\n// references stay valid).\npageBreakParagraph.delete();\nblankParagraph2.delete();\nblankParagraph1.delete();\nawait context.sync();\n", "ps1": "# The \"MBS Datasets\" Heading-3 paragraph loses its paragraph style / spacing\n# and instead carries the equivalent look (font, color, size) as direct\n# formatting on the paragraph mark, while the trailing blank paragraphs and\n# the paragraph that used to hold the page-break run are collapsed into it\n# (the page-break run itself is kept, now ending the \"MBS Datasets\"\n# paragraph).\n$d = $word.ActiveDocument\n\n# Locate the first \"MBS Datasets\" Heading 3 paragraph that is immediately\n# followed by two empty paragraphs and a page-break paragraph (the second\n# \"MBS Datasets\" heading later in the document must stay untouched).\n$count = $d.Paragraphs.Count\n$targetIndex = -1\nfor ($i = 1; $i -le $count - 3; $i++) {\n    $t0 = $d.Paragraphs.Item($i).Range.Text\n    $t1 = $d.Paragraphs.Item($i + 1).Range.Text\n    $t2 = $d.Paragraphs.Item($i + 2).Range.Text\n    $t3 = $d.Paragraphs.Item($i + 3).Range.Text\n    if ($t0 -eq \"MBS Datasets`r\" -and $t1 -eq \"`r\" -and $t2 -eq \"`r\" -and $t3 -eq ([char]12 + \"`r\")) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not locate the target 'MBS Datasets' heading paragraph.\"\n}\n\n# Replace the heading paragraph's OOXML: drop pStyle/spacing, bake the\n# Heading-3 look into the paragraph mark's rPr, keep the two existing runs\n# untouched, and append the page-break run that used to live in its own\n# paragraph.\n$headingRange = $d.Paragraphs.Item($targetIndex).Range\n$newParagraphXml = '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n'<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n'<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n'<pkg:xmlData>' +\n'<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n'<w:body>' +\n'<w:p w14:paraId=\"708C0CEA\" w14:textId=\"38B3AB9C\" w:rsidR=\"007325B7\" w:rsidRPr=\"000E1730\" w:rsidRDefault=\"007325B7\" w:rsidP=\"007325B7\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n'<w:pPr>' +\n'<w:rPr>' +\n'<w:rFonts w:eastAsiaTheme=\"majorEastAsia\" w:cstheme=\"majorBidi\"/>' +\n'<w:color w:val=\"0F4761\" w:themeColor=\"accent1\" w:themeShade=\"BF\"/>' +\n'<w:sz w:val=\"28\"/>' +\n'<w:szCs w:val=\"28\"/>' +\n'<w:u w:val=\"single\"/>' +\n'</w:rPr>' +\n'</w:pPr>' +\n'<w:r>' +\n'<w:rPr>' +\n'<w:u w:val=\"single\"/>' +\n'</w:rPr>' +\n'<w:t>MBS</w:t>' +\n'</w:r>' +\n'<w:r w:rsidRPr=\"000E1730\">' +\n'<w:rPr>' +\n'<w:u w:val=\"single\"/>' +\n'</w:rPr>' +\n'<w:t xml:space=\"preserve\"> Datasets</w:t>' +\n'</w:r>' +\n'<w:r>' +\n'<w:rPr>' +\n'<w:b/>' +\n'<w:bCs/>' +\n'</w:rPr>' +\n'<w:br w:type=\"page\"/>' +\n'</w:r>' +\n'</w:p>' +\n'</w:body>' +\n'</w:document>' +\n'</pkg:xmlData>' +\n'</pkg:part>' +\n'</pkg:package>'\n\n$headingRange.InsertXML($newParagraphXml)\n\n# Remove the two blank paragraphs and the now-redundant paragraph that used\n# to hold only the page-break run (delete from the bottom up so the other\n# indices stay valid).\n$d.Paragraphs.Item($targetIndex + 3).Range.Delete()\n$d.Paragraphs.Item($targetIndex + 2).Range.Delete()\n$d.Paragraphs.Item($targetIndex + 1).Range.Delete()\n"}
